$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{
        G = 236.082283; H = 472.164566; I = 0.356044175747874; J = 0.2923054205394457;
        M = 0.4333629999999999; N = 0.8667259999999999;
        Q = 102.309326407729; R = 409.237305630916; S = 0.356044175747874; T = 0.2923054205394457
    }
    3 = @{
        I = 0.1772841109604352; J = 0.2183202681257223;
        M = 0.4333629999999999; N = 0.8667259999999999;
        Q = 50.94260547039299; R = 305.655632822358; S = 0.1772841109604352; T = 0.2183202681257223
    }
    4 = @{
        G = 65.58120733333334; H = 196.743622; I = 0.09890537575641412; J = 0.1217991168934175;
        M = 0.4333629999999999; N = 0.8667259999999999;
        Q = 28.42046875359533; R = 170.522812521572; S = 0.09890537575641412; T = 0.1217991168934175
    }
    5 = @{
        G = 137.8159335; H = 275.631867; I = 0.2078451623916705; J = 0.1706368808656166;
        M = 0.4333629999999999; N = 0.8667259999999999;
        Q = 59.72432638936049; R = 238.897305557442; S = 0.2078451623916705; T = 0.1706368808656166
    }
    6 = @{
        G = 39.78016766666667; H = 119.340503; I = 0.05999390055030331; J = 0.07388075774581522;
        M = 0.4333629999999999; N = 0.8667259999999999;
        Q = 17.23925280052967; R = 103.435516803178; S = 0.05999390055030331; T = 0.07388075774581522
    }
    7 = @{
        G = 66.258798; H = 198.776394; I = 0.09992727459330304; J = 0.1230575558299827;
        M = 0.4333629999999999; N = 0.8667259999999999;
        Q = 28.714111477674; R = 172.284668866044; S = 0.09992727459330304; T = 0.1230575558299827
    }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}

$wb.Save()
